$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.444.66"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "2.685.83"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "521.77"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "146.20"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "2.702.89"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "3.157.51"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "60.442.42"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "21.27"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.761.24"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000139"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "351.07"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("D20").Value = "4.54"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "6.33"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "62.88"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +4.76%  "
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "0.0₃0815"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").Value = "6.87"
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("D33").Value = "19.06"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "148.37"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  +7.09%  "
$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "0.949"
$ws.Range("E36").Value = "  -6.96%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +7.36%  "
$ws.Range("E38").Value = "  +10.69%  "
$ws.Range("D39").Value = "0.878"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").Value = "36.81"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").Value = "282.64"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.613"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.96"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "0.0987"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "2.128.48"
$ws.Range("E47").Value = "  +6.92%  "
$ws.Range("D48").Value = "0.0539"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").Value = "4.87"
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("D50").Value = "0.0235"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "10.45"
$ws.Range("E51").Value = "  +1.82%  "
